$d = $word.ActiveDocument

$replacements = @(
    @("762÷7=", "984÷8="),
    @("143÷6=", "848÷8="),
    @("145÷7=", "580÷6="),
    @("479÷5=", "308÷7="),
    @("616÷8=", "829÷7="),
    @("589÷3=", "223÷2="),
    @("992÷6=", "573÷3="),
    @("941÷7=", "862÷2="),
    @("812÷7=", "370÷9="),
    @("283÷5=", "182÷2="),
    @("561÷5=", "723÷4="),
    @("456÷3=", "567÷8="),
    @("857÷6=", "250÷4="),
    @("176÷3=", "422÷6="),
    @("455÷3=", "775÷5="),
    @("582÷8=", "537÷9="),
    @("477÷4=", "250÷5="),
    @("754÷5=", "247÷7="),
    @("947÷9=", "568÷7="),
    @("177÷8=", "352÷3="),
    @("679÷2=", "494÷4="),
    @("157÷5=", "518÷9="),
    @("138÷9=", "148÷8="),
    @("551÷2=", "118÷4="),
    @("792÷3=", "200÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
